# "more skill data balancing"
#
# Changes applied to the "Skills" worksheet:
#  1. Row 4, column I ("Skill4" for the Improvise/Tinker entry) changes from
#     "STUDY \xF0\x9F\x92\xA1\xF0\x9F\x92\xA1" to a brand new skill roll value
#     "YANK! \xF0\x9F\x92\xA1\xF0\x9F\x92\xA1\xF0\x9F\x92\xA1\xF0\x9F\x94\x8A\xF0\x9F\x94\x8A".
#     (The three obsolete "FIGHT!", "SPRINT!" and "Smash and Grab" shared
#     strings disappear naturally because nothing references them anymore.)
#  2. The whole "Qty" (E) column for the lower reference table (rows 8-31) is
#     removed entirely -- not just blanked -- so the cells no longer appear
#     in the sheet at all.
#  3. The small scratch table that lived in rows 34-36 (Smash and Grab /
#     Script Kiddie / Hacktivist notes) is deleted completely.
#  4. The active selection moves from K6 to E9:E10, reflecting the area the
#     author was last working in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skills")
$ws.Activate()

# 1. Update the single changed skill-roll cell.
$ws.Range("I4").Value = "YANK! 💡💡💡🔊🔊"

# 2. Remove the Qty (E) values for the lower table entirely (cells vanish,
#    not merely cleared of their value).
$ws.Range("E8:E31").Clear()

# 3. Delete the leftover scratch rows at the bottom of the sheet.
$ws.Rows("34:36").Delete()

# 4. Leave the selection on the range the author ended up editing.
$ws.Range("E9:E10").Select()
